$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 92, pushing the existing row 92 (and all rows
# below it) down by one. This grows the used range from A1:R160 to A1:R161.
$ws.Rows.Item(92).Insert()

# Populate the newly inserted row 92 with the new weekly price record.
$ws.Range("A92").Value = 8
$ws.Range("B92").Value = 'Terminal La Palmera de La Serena'
$ws.Range("C92").Value = 'Coquimbo'
$ws.Range("D92").Value = 44978
$ws.Range("E92").Value = 4
$ws.Range("F92").Value = 100112052
$ws.Range("G92").Value = 'Albahaca'
$ws.Range("H92").Value = 'Sin especificar'
$ws.Range("I92").Value = 'Primera'
$ws.Range("J92").Value = 500
$ws.Range("K92").Value = 4500
$ws.Range("L92").Value = 5000
$ws.Range("M92").Value = 4750
$ws.Range("N92").Value = '$/docena de matas'
$ws.Range("O92").Value = 'Provincia del Elquí'
$ws.Range("P92").Value = 792
$ws.Range("Q92").Value = 6
$ws.Range("R92").Value = 'Hortaliza'
